$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8147.6665
$ws.Range("I62").Value = 4999.5
$ws.Range("K62").Value = 4999.5
$ws.Range("M62").Value = -4375.5
$ws.Range("H65").Value = 8147.6665
$ws.Range("I65").Value = 4999.5
$ws.Range("K65").Value = 24997.5
$ws.Range("M65").Value = -21877.5
$ws.Range("H86").Value = 2615.7058
$ws.Range("I86").Value = 3642.75
$ws.Range("J86").Value = 1702.7778
$ws.Range("K86").Value = 3642.75
$ws.Range("L86").Value = 1702.7778
$ws.Range("M86").Value = -2519.75
$ws.Range("N86").Value = -3948.7778
$ws.Range("H89").Value = 2615.7058
$ws.Range("I89").Value = 3642.75
$ws.Range("J89").Value = 1702.7778
$ws.Range("K89").Value = 18213.75
$ws.Range("L89").Value = 8513.889000000001
$ws.Range("M89").Value = -12597.75
$ws.Range("N89").Value = -19745.889
$ws.Range("H132").Value = 4529.0435
$ws.Range("I132").Value = 3740.4285
$ws.Range("K132").Value = 11221.2855
$ws.Range("M132").Value = -8691.2855
$ws.Range("H137").Value = 2499.889
$ws.Range("I137").Value = 2249.8333
$ws.Range("K137").Value = 6749.499899999999
$ws.Range("M137").Value = -4199.499899999999
$ws.Range("H141").Value = 1528.8462
$ws.Range("I141").Value = 1572.9166
$ws.Range("K141").Value = 4718.7498
$ws.Range("M141").Value = 461.2502000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9673.666999999999
$ws.Range("J2").Value = 10010.5
$ws.Range("L2").Value = 10010.5
$ws.Range("N2").Value = -10236.5
$ws.Range("H74").Value = 1619.3
$ws.Range("I74").Value = 1688.3334
$ws.Range("J74").Value = 998
$ws.Range("K74").Value = 1688.3334
$ws.Range("L74").Value = 998
$ws.Range("M74").Value = -814.3334
$ws.Range("N74").Value = -2746
$ws.Range("H77").Value = 1619.3
$ws.Range("I77").Value = 1688.3334
$ws.Range("J77").Value = 998
$ws.Range("K77").Value = 8441.666999999999
$ws.Range("L77").Value = 4990
$ws.Range("M77").Value = -4073.666999999999
$ws.Range("N77").Value = -13726
$ws.Range("H102").Value = 1998.2858
$ws.Range("I102").Value = 1164.6666
$ws.Range("K102").Value = 1164.6666
$ws.Range("M102").Value = 457.3334
$ws.Range("H116").Value = 9673.666999999999
$ws.Range("J116").Value = 10010.5
$ws.Range("L116").Value = 10010.5
$ws.Range("N116").Value = -14598.5
$ws.Range("H122").Value = 4999
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9673.666999999999
$ws.Range("J3").Value = 10010.5
$ws.Range("L3").Value = 10010.5
$ws.Range("N3").Value = -10238.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2120.875
$ws.Range("I31").Value = 1742.5
$ws.Range("K31").Value = 1742.5
$ws.Range("M31").Value = -1447.5
$ws.Range("H34").Value = 2120.875
$ws.Range("I34").Value = 1742.5
$ws.Range("K34").Value = 1742.5
$ws.Range("M34").Value = -1540.5
$ws.Range("H132").Value = 5089.5
$ws.Range("I132").Value = 5837
$ws.Range("J132").Value = 2099.5
$ws.Range("K132").Value = 17511
$ws.Range("L132").Value = 6298.5
$ws.Range("M132").Value = -14981
$ws.Range("N132").Value = -11358.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 51.714287
$ws.Range("I14").Value = 51.714287
$ws.Range("K14").Value = 155.142861
$ws.Range("M14").Value = 17.85713900000002
$ws.Range("H57").Value = 200
$ws.Range("I57").Value = 200
$ws.Range("K57").Value = 600
$ws.Range("M57").Value = -41
$ws.Range("H80").Value = 11874.75
$ws.Range("J80").Value = 14999.667
$ws.Range("L80").Value = 44999.001
$ws.Range("N80").Value = -46871.001
$ws.Range("H83").Value = 11874.75
$ws.Range("J83").Value = 14999.667
$ws.Range("L83").Value = 134997.003
$ws.Range("N83").Value = -144357.003
$ws.Range("H139").Value = 1677.6666
$ws.Range("I139").Value = 2500
$ws.Range("J139").Value = 33
$ws.Range("K139").Value = 7500
$ws.Range("L139").Value = 99
$ws.Range("M139").Value = -2360
$ws.Range("N139").Value = -10379

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1947
$ws.Range("I7").Value = 1947
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1947
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1835
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 2142.4285
$ws.Range("I22").Value = 1166.1666
$ws.Range("K22").Value = 1166.1666
$ws.Range("M22").Value = -871.1666
$ws.Range("H27").Value = 2142.4285
$ws.Range("I27").Value = 1166.1666
$ws.Range("K27").Value = 1166.1666
$ws.Range("M27").Value = -1059.1666
$ws.Range("H46").Value = 1739.4286
$ws.Range("I46").Value = 1313.3334
$ws.Range("J46").Value = 2231.077
$ws.Range("K46").Value = 1313.3334
$ws.Range("L46").Value = 2231.077
$ws.Range("M46").Value = -1125.3334
$ws.Range("N46").Value = -2607.077
$ws.Range("H126").Value = 1947
$ws.Range("I126").Value = 1947
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5841
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3371
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4270.0713
$ws.Range("I132").Value = 4628.4
$ws.Range("K132").Value = 13885.2
$ws.Range("M132").Value = -11355.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4750
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050
$ws.Range("H126").Value = 1581.9286
$ws.Range("I126").Value = 1599.7693
$ws.Range("K126").Value = 4799.3079
$ws.Range("M126").Value = -2329.3079
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 1800.75
$ws.Range("I136").Value = 914.4783
$ws.Range("J136").Value = 5877.6
$ws.Range("K136").Value = 2743.4349
$ws.Range("L136").Value = 17632.8
$ws.Range("M136").Value = -193.4349000000002
$ws.Range("N136").Value = -22732.8
